$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 284
$ws.Cells.Item(284,1).Value = "●"
$ws.Cells.Item(284,2).Value = "'"
$ws.Cells.Item(284,3).Value = "'"
$ws.Cells.Item(284,4).Value = "'2863"
$ws.Cells.Item(284,5).Value = "Event year"
$ws.Cells.Item(284,6).Value = "1: 1933"
$ws.Cells.Item(284,7).Value = "1: 1936"
$ws.Cells.Item(284,8).Value = 0
$ws.Cells.Item(284,9).Value = "'2000"
$ws.Cells.Item(284,10).Value = 4
$ws.Cells.Item(284,11).Value = 0.011254
$ws.Cells.Item(284,12).Value = "Sonia"
$ws.Cells.Item(284,13).Value = "11/14/18 13:12:00"
$ws.Range("A283:M283").Copy()
$ws.Range("A284:M284").PasteSpecial(-4122)
$ws.Rows.Item(284).RowHeight = 16

# Row 285
$ws.Cells.Item(285,1).Value = "●"
$ws.Cells.Item(285,2).Value = "'"
$ws.Cells.Item(285,3).Value = "'"
$ws.Cells.Item(285,4).Value = "'2863"
$ws.Cells.Item(285,5).Value = "Event year"
$ws.Cells.Item(285,6).Value = "1: 1942"
$ws.Cells.Item(285,7).Value = "1: 1945"
$ws.Cells.Item(285,8).Value = 0
$ws.Cells.Item(285,9).Value = "'2013"
$ws.Cells.Item(285,10).Value = 4
$ws.Cells.Item(285,11).Value = 0.011254
$ws.Cells.Item(285,12).Value = "Sonia"
$ws.Cells.Item(285,13).Value = "11/14/18 13:12:00"
$ws.Range("A283:M283").Copy()
$ws.Range("A285:M285").PasteSpecial(-4122)
$ws.Rows.Item(285).RowHeight = 16

# Row 286
$ws.Cells.Item(286,1).Value = "●"
$ws.Cells.Item(286,2).Value = "'"
$ws.Cells.Item(286,3).Value = "'"
$ws.Cells.Item(286,4).Value = "'5362"
$ws.Cells.Item(286,5).Value = "Event month"
$ws.Cells.Item(286,6).Value = "1: 1890"
$ws.Cells.Item(286,7).Value = "1: 1893"
$ws.Cells.Item(286,8).Value = 0
$ws.Cells.Item(286,9).Value = "June"
$ws.Cells.Item(286,10).Value = 4
$ws.Cells.Item(286,11).Value = 0.011473
$ws.Cells.Item(286,12).Value = "Sonia"
$ws.Cells.Item(286,13).Value = "11/14/18 13:17:00"
$ws.Range("A283:M283").Copy()
$ws.Range("A286:M286").PasteSpecial(-4122)
$ws.Rows.Item(286).RowHeight = 16

# Row 287
$ws.Cells.Item(287,1).Value = "●"
$ws.Cells.Item(287,2).Value = "'"
$ws.Cells.Item(287,3).Value = "'"
$ws.Cells.Item(287,4).Value = "'5362"
$ws.Cells.Item(287,5).Value = "Event month"
$ws.Cells.Item(287,6).Value = "1: 1903"
$ws.Cells.Item(287,7).Value = "1: 1909"
$ws.Cells.Item(287,8).Value = 0
$ws.Cells.Item(287,9).Value = "January"
$ws.Cells.Item(287,10).Value = 7
$ws.Cells.Item(287,11).Value = 0.020077
$ws.Cells.Item(287,12).Value = "Sonia"
$ws.Cells.Item(287,13).Value = "11/14/18 13:17:00"
$ws.Range("A283:M283").Copy()
$ws.Range("A287:M287").PasteSpecial(-4122)
$ws.Rows.Item(287).RowHeight = 16

# Row 288
$ws.Cells.Item(288,1).Value = "●"
$ws.Cells.Item(288,2).Value = "'"
$ws.Cells.Item(288,3).Value = "'"
$ws.Cells.Item(288,4).Value = "'5362"
$ws.Cells.Item(288,5).Value = "Event year"
$ws.Cells.Item(288,6).Value = "1: 1895"
$ws.Cells.Item(288,7).Value = "1: 1898"
$ws.Cells.Item(288,8).Value = 0
$ws.Cells.Item(288,9).Value = "'2009"
$ws.Cells.Item(288,10).Value = 4
$ws.Cells.Item(288,11).Value = 0.011473
$ws.Cells.Item(288,12).Value = "Sonia"
$ws.Cells.Item(288,13).Value = "11/14/18 13:17:00"
$ws.Range("A283:M283").Copy()
$ws.Range("A288:M288").PasteSpecial(-4122)
$ws.Rows.Item(288).RowHeight = 16

# Row 289
$ws.Cells.Item(289,1).Value = "●"
$ws.Cells.Item(289,2).Value = "'"
$ws.Cells.Item(289,3).Value = "'"
$ws.Cells.Item(289,4).Value = "'5362"
$ws.Cells.Item(289,5).Value = "Event year"
$ws.Cells.Item(289,6).Value = "1: 1911"
$ws.Cells.Item(289,7).Value = "1: 1914"
$ws.Cells.Item(289,8).Value = 0
$ws.Cells.Item(289,9).Value = "'2010"
$ws.Cells.Item(289,10).Value = 4
$ws.Cells.Item(289,11).Value = 0.011473
$ws.Cells.Item(289,12).Value = "Sonia"
$ws.Cells.Item(289,13).Value = "11/14/18 13:18:00"
$ws.Range("A283:M283").Copy()
$ws.Range("A289:M289").PasteSpecial(-4122)
$ws.Rows.Item(289).RowHeight = 16

# Row 290
$ws.Cells.Item(290,1).Value = "●"
$ws.Cells.Item(290,2).Value = "'"
$ws.Cells.Item(290,3).Value = "'"
$ws.Cells.Item(290,4).Value = "'5362"
$ws.Cells.Item(290,5).Value = "B"
$ws.Cells.Item(290,6).Value = "1: 1911"
$ws.Cells.Item(290,7).Value = "1: 1914"
$ws.Cells.Item(290,8).Value = 0
$ws.Cells.Item(290,9).Value = "'2010"
$ws.Cells.Item(290,10).Value = 4
$ws.Cells.Item(290,11).Value = 0.011473
$ws.Cells.Item(290,12).Value = "Sonia"
$ws.Cells.Item(290,13).Value = "11/14/18 13:18:00"
$ws.Range("A283:M283").Copy()
$ws.Range("A290:M290").PasteSpecial(-4122)
$ws.Rows.Item(290).RowHeight = 16

# Row 291
$ws.Cells.Item(291,1).Value = "●"
$ws.Cells.Item(291,2).Value = "'"
$ws.Cells.Item(291,3).Value = "'"
$ws.Cells.Item(291,4).Value = "'5362"
$ws.Cells.Item(291,5).Value = "B"
$ws.Cells.Item(291,6).Value = "1: 1903"
$ws.Cells.Item(291,7).Value = "1: 1909"
$ws.Cells.Item(291,8).Value = 0
$ws.Cells.Item(291,9).Value = "January"
$ws.Cells.Item(291,10).Value = 7
$ws.Cells.Item(291,11).Value = 0.020077
$ws.Cells.Item(291,12).Value = "Sonia"
$ws.Cells.Item(291,13).Value = "11/14/18 13:18:00"
$ws.Range("A283:M283").Copy()
$ws.Range("A291:M291").PasteSpecial(-4122)
$ws.Rows.Item(291).RowHeight = 16

# Row 292
$ws.Cells.Item(292,1).Value = "●"
$ws.Cells.Item(292,2).Value = "'"
$ws.Cells.Item(292,3).Value = "'"
$ws.Cells.Item(292,4).Value = "'5362"
$ws.Cells.Item(292,5).Value = "A"
$ws.Cells.Item(292,6).Value = "1: 1895"
$ws.Cells.Item(292,7).Value = "1: 1898"
$ws.Cells.Item(292,8).Value = 0
$ws.Cells.Item(292,9).Value = "'2009"
$ws.Cells.Item(292,10).Value = 4
$ws.Cells.Item(292,11).Value = 0.011473
$ws.Cells.Item(292,12).Value = "Sonia"
$ws.Cells.Item(292,13).Value = "11/14/18 13:18:00"
$ws.Range("A283:M283").Copy()
$ws.Range("A292:M292").PasteSpecial(-4122)
$ws.Rows.Item(292).RowHeight = 16

# Row 293
$ws.Cells.Item(293,1).Value = "●"
$ws.Cells.Item(293,2).Value = "'"
$ws.Cells.Item(293,3).Value = "'"
$ws.Cells.Item(293,4).Value = "'5362"
$ws.Cells.Item(293,5).Value = "A"
$ws.Cells.Item(293,6).Value = "1: 1890"
$ws.Cells.Item(293,7).Value = "1: 1893"
$ws.Cells.Item(293,8).Value = 0
$ws.Cells.Item(293,9).Value = "June"
$ws.Cells.Item(293,10).Value = 4
$ws.Cells.Item(293,11).Value = 0.011473
$ws.Cells.Item(293,12).Value = "Sonia"
$ws.Cells.Item(293,13).Value = "11/14/18 13:17:00"
$ws.Range("A283:M283").Copy()
$ws.Range("A293:M293").PasteSpecial(-4122)
$ws.Rows.Item(293).RowHeight = 16
